# Lab Notebooks -> "Final presentations for NOAA" edit
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Handout master / notes master: update the fixed "date" field text
# ---------------------------------------------------------------------------
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "7/28/2023"

$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "7/28/2023"

# ---------------------------------------------------------------------------
# 2. Slide 1 (title slide): speaker name / pronoun / affiliation / event
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# 2a. Speaker name: "Anshu" + " Dubey" -> "David E. Bernholdt"
$nameShape = $s1.Shapes.Item(3)
$nameRange = $nameShape.TextFrame.TextRange
$firstWord = $nameRange.Characters(1, 5)   # "Anshu"
$firstWord.Text = ""
$nameRange.Runs(1).Text = "David E. Bernholdt"
# widen the box to fit the new (longer) name
$nameShape.Width = 216.74623047244094      # -> 2752677 EMU

# 2b. Pronoun: "(she/her)" -> "(he/him)"
$pronounShape = $s1.Shapes.Item(4)
$pronounShape.TextFrame.TextRange.Text = "(he/him)"
# shift box right to line up with the widened name box
$pronounShape.Left = 465.7102862204724     # -> 5914520 EMU

# 2c. Affiliation: "Argonne National Laboratory" -> "Oak Ridge National Laboratory"
$labShape = $s1.Shapes.Item(5)
$labShape.TextFrame.TextRange.Text = "Oak Ridge National Laboratory"

# 2d. Event: "Better Scientific Software tutorial @ ISC23" -> "... @ NOAA Global Systems Laboratory "
$eventShape = $s1.Shapes.Item(6)
$eventShape.TextFrame.TextRange.Text = "Better Scientific Software tutorial @ NOAA Global Systems Laboratory "

# ---------------------------------------------------------------------------
# 3. Slide 2 (license / citation slide): update citation paragraph
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$contentShape = $s2.Shapes.Item(2)
$citation = $contentShape.TextFrame.TextRange.Paragraphs(3)

# Runs in this paragraph (before edit):
#  1: "The requested citation the overall tutorial is: "                 (unchanged)
#  2: "Anshu Dubey and David M. Rogers, ... DOI: "                       (text + font)
#  3: "10.6084/m9.figshare" + "22790762" (two merged sub-runs)           (text + font + color + underline)
#  4: "."                                                                (font)

# 3a. Remove the duplicate/second DOI sub-run ".22790762" first (same
#     formatting as the visible merged run, so indices stay stable),
#     then rewrite the remaining run text -- avoids leaving stray text.
$dupRun = $citation.Characters(203, 9)      # ".22790762"
$dupRun.Text = ""

# 3b. Rewrite the long citation sentence + switch theme font to major/heading
$r2 = $citation.Runs(2)
$r2.Text = "David E. Bernholdt, Anshu Dubey, and Patricia A. Grubel, Better Scientific Software tutorial, in NOAA Global Systems Laboratory, Boulder, Colorado, 2023. DOI: "
$r2.Font.Name = "+mj-lt"

# 3c. Rewrite the DOI run, restyle it (no underline, no strike, blue link color)
$r3 = $citation.Runs(3)
$r3.Text = "10.6084/m9.figshare.23796606"
$r3.Font.Name = "+mj-lt"
$r3.Font.Underline = $false
$r3.Font.Strike = 0
$r3.Font.Color.RGB = 0xE27A2A

# 3d. Trailing "." run: only the theme font changes
$r4 = $citation.Runs(4)
$r4.Font.Name = "+mj-lt"

Write-Host "edit.ps1 complete"
